# Fruta / hortaliza, semanal
# Insert a new weekly record at row 575 (pushing the existing rows 575-601
# down to 576-602) and populate it with this week's data for
# "Terminal La Palmera de La Serena" - Choclo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 575:601 down by one to make room for the new record.
$ws.Rows.Item(575).Insert()

# Fill in the new row 575 with the new weekly observation.
$ws.Cells.Item(575, 1).Value = 8
$ws.Cells.Item(575, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(575, 3).Value = "Coquimbo"
$ws.Cells.Item(575, 4).Value = (Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(575, 5).Value = 4
$ws.Cells.Item(575, 6).Value = 100112024
$ws.Cells.Item(575, 7).Value = "Choclo"
$ws.Cells.Item(575, 8).Value = "Dulce o Americano"
$ws.Cells.Item(575, 9).Value = "Primera"
$ws.Cells.Item(575, 10).Value = 25000
$ws.Cells.Item(575, 11).Value = 450
$ws.Cells.Item(575, 12).Value = 500
$ws.Cells.Item(575, 13).Value = 475
$ws.Cells.Item(575, 14).Value = "$/unidad"
$ws.Cells.Item(575, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(575, 16).Value = 475
$ws.Cells.Item(575, 17).Value = 1
$ws.Cells.Item(575, 18).Value = "Hortaliza"
